$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) label updates ---
$ws.Range("A1").Value = "tag"
$ws.Range("D1").Value = "formula"
$ws.Range("F1").Value = "attr"

# --- Row 2: update measured value ---
$ws.Range("B2").Value = 3.15

# --- Row 3: update measured value ---
$ws.Range("B3").Value = 3.15

# --- Row 4: update measured value and formula branch labels ---
$ws.Range("B4").Value = 3.15
$ws.Range("D4").Formula = '=IF(B4>C4,"red","yellow")'

# --- Row 5: convert text "true" values to numeric, and update formula branch labels ---
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Formula = '=IF(B5=C5,"green","blue")'

# --- Drop the (invisible, borderId=0) leftover border flag on the data cells in
#     columns A, D and F so they fall back to the plain default style, matching
#     the cleaned-up style table produced when the workbook was re-saved ---
$ws.Range("A2:A5").Borders.LineStyle = 0
$ws.Range("D2:D5").Borders.LineStyle = 0
$ws.Range("F2:F5").Borders.LineStyle = 0
